$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings keep their exact
# formatting (trailing zeros, leading zeros, multi-dot separators) instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "31.224.06"
$ws.Cells.Item(2, 5).Value = "  +2.81%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.997.50"
$ws.Cells.Item(3, 5).Value = "  +6.72%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "0.9993"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "0.7818"
$ws.Cells.Item(5, 5).Value = "  +65.50%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "257.10"
$ws.Cells.Item(6, 5).Value = "  +5.07%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.9998"
$ws.Cells.Item(7, 5).Value = "  -0.05%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.3525"
$ws.Cells.Item(8, 5).Value = "  +22.94%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "28.71"
$ws.Cells.Item(9, 5).Value = "  +31.88%  "

# Row 10
$ws.Cells.Item(10, 2).Value = "OKB"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(10, 4).Value = "44.38"
$ws.Cells.Item(10, 5).Value = "  +3.25%  "

# Row 11
$ws.Cells.Item(11, 2).Value = "Dogecoin"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(11, 4).Value = "0.07040"
$ws.Cells.Item(11, 5).Value = "  +8.56%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "Polygon"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(12, 4).Value = "0.8573"
$ws.Cells.Item(12, 5).Value = "  +17.86%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "TRON"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(13, 4).Value = "0.08204"
$ws.Cells.Item(13, 5).Value = "  +5.22%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "101.35"
$ws.Cells.Item(14, 5).Value = "  +1.23%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "WrappedEther"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(15, 4).Value = "1.999.27"
$ws.Cells.Item(15, 5).Value = "  +6.85%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "Polkadot"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(16, 4).Value = "5.584"
$ws.Cells.Item(16, 5).Value = "  +8.13%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "Avalanche"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(17, 4).Value = "15.42"
$ws.Cells.Item(17, 5).Value = "  +17.87%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "BitcoinCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(18, 4).Value = "274.31"
$ws.Cells.Item(18, 5).Value = "  -3.27%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "WrappedBTC"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(19, 4).Value = "31.229.60"
$ws.Cells.Item(19, 5).Value = "  +2.89%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(20, 4).Value = "5.964"
$ws.Cells.Item(20, 5).Value = "  +11.96%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "ShibaInu"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(21, 4).Value = "0.000007953"
$ws.Cells.Item(21, 5).Value = "  +6.45%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(22, 4).Value = "2.263.04"
$ws.Cells.Item(22, 5).Value = "  +7.10%  "

# Row 23
$ws.Cells.Item(23, 2).Value = "Dai"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(23, 4).Value = "0.9999"
$ws.Cells.Item(23, 5).Value = "  -0.04%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "BinanceUSD"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(24, 4).Value = "0.9988"
$ws.Cells.Item(24, 5).Value = "  -0.11%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "Chainlink"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(25, 4).Value = "7.140"
$ws.Cells.Item(25, 5).Value = "  +12.90%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "Cosmos"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(26, 4).Value = "10.05"
$ws.Cells.Item(26, 5).Value = "  +11.37%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "0.1485"
$ws.Cells.Item(27, 5).Value = "  +53.75%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "Monero"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(28, 4).Value = "163.79"
$ws.Cells.Item(28, 5).Value = "  +0.56%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "EthereumClassic"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(29, 4).Value = "19.97"
$ws.Cells.Item(29, 5).Value = "  +5.53%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "LidoDAOToken"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(30, 4).Value = "2.360"
$ws.Cells.Item(30, 5).Value = "  +24.84%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).Value = "1.608"
$ws.Cells.Item(31, 5).Value = "  +7.97%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "4.621"
$ws.Cells.Item(32, 5).Value = "  +9.41%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Toncoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(33, 4).Value = "1.358"
$ws.Cells.Item(33, 5).Value = "  +2.69%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(34, 4).Value = "4.428"
$ws.Cells.Item(34, 5).Value = "  +6.92%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Hedera"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(35, 4).Value = "0.05210"
$ws.Cells.Item(35, 5).Value = "  +8.48%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).Value = "1.228"
$ws.Cells.Item(36, 5).Value = "  +9.24%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "ImmutableX"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(37, 4).Value = "0.7747"
$ws.Cells.Item(37, 5).Value = "  +12.52%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "HuobiToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(38, 4).Value = "2.812"
$ws.Cells.Item(38, 5).Value = "  +3.36%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(39, 4).Value = "0.02004"
$ws.Cells.Item(39, 5).Value = "  +5.60%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "MXToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(40, 4).Value = "2.912"
$ws.Cells.Item(40, 5).Value = "  +2.62%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(41, 4).Value = "6.710"
$ws.Cells.Item(41, 5).Value = "  +6.81%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Aave"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(42, 4).Value = "79.50"
$ws.Cells.Item(42, 5).Value = "  +4.05%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "TheSandbox"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(43, 4).Value = "0.4721"
$ws.Cells.Item(43, 5).Value = "  +12.11%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "RenderToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(44, 4).Value = "2.154"
$ws.Cells.Item(44, 5).Value = "  +10.42%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Quant"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(45, 4).Value = "106.61"
$ws.Cells.Item(45, 5).Value = "  +5.67%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "TrustWalletToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(46, 4).Value = "0.8546"
$ws.Cells.Item(46, 5).Value = "  +3.90%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "PaxDollar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(47, 4).Value = "1.000"
$ws.Cells.Item(47, 5).Value = "  +0.07%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Aptos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(48, 4).Value = "7.729"
$ws.Cells.Item(48, 5).Value = "  +10.25%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value = "9.951"
$ws.Cells.Item(49, 5).Value = "  +2.43%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Decentraland"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(50, 4).Value = "0.4325"
$ws.Cells.Item(50, 5).Value = "  +10.83%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Elrond"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(51, 4).Value = "36.71"
$ws.Cells.Item(51, 5).Value = "  +5.27%  "
